# Insert two new columns before column G ("SECTION"), shifting the
# existing G:U header columns right to I:W.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G:H").Insert()

# Copy the header formatting from the neighbouring "ADMIN ASSIGNOR" header
# cell (F4) onto the two newly inserted header cells, then set their text.
$ws.Range("F4").Copy()
$ws.Range("G4:H4").PasteSpecial(-4122)
$ws.Range("G4").Value = "SUBJECT"
$ws.Range("H4").Value = "DESCRIPTION"

# Match the authored column widths as closely as the host allows.
$ws.Columns.Item(7).ColumnWidth = 32.5
$ws.Columns.Item(8).ColumnWidth = 28.6

# Restore the view/selection state recorded in the edited workbook.
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q19").Select()

# Window position recorded against the workbook view.
$excel.ActiveWindow.Left = 6720
$excel.ActiveWindow.Top = 1590
